$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "62.522.32"; E = "  -0.94%  " },
    @{ Row = 3; D = "3.435.95"; E = "  -1.47%  " },
    @{ Row = 4; D = "'0.999"; E = "  -0.15%  " },
    @{ Row = 5; D = "'577.89"; E = "  -1.09%  " },
    @{ Row = 6; D = "'147.26"; E = "  -0.54%  " },
    @{ Row = 7; D = $null; E = "  -0.02%  " },
    @{ Row = 8; D = $null; E = "  +0.25%  " },
    @{ Row = 9; D = "'7.94"; E = "  +3.27%  " },
    @{ Row = 10; D = $null; E = "  -2.09%  " },
    @{ Row = 11; D = $null; E = "  +2.40%  " },
    @{ Row = 12; D = "4.026.67"; E = "  -1.57%  " },
    @{ Row = 13; D = $null; E = "  +2.29%  " },
    @{ Row = 14; D = "'28.19"; E = "  -5.48%  " },
    @{ Row = 15; D = "3.445.35"; E = "  -1.32%  " },
    @{ Row = 16; D = $null; E = "  -0.85%  " },
    @{ Row = 17; D = "62.561.09"; E = "  -1.05%  " },
    @{ Row = 18; D = "'6.36"; E = "  +0.53%  " },
    @{ Row = 19; D = "'14.53"; E = "  +1.41%  " },
    @{ Row = 20; D = "'9.05"; E = "  -3.14%  " },
    @{ Row = 21; D = "'385.98"; E = "  -0.98%  " },
    @{ Row = 22; D = "'75.08"; E = "  +0.05%  " },
    @{ Row = 23; D = "'0.559"; E = "  -0.83%  " },
    @{ Row = 24; D = $null; E = "  +0.07%  " },
    @{ Row = 25; D = "3.582.94"; E = "  -1.36%  " },
    @{ Row = 26; D = "'0.0000114"; E = "  -2.49%  " },
    @{ Row = 27; D = $null; E = "  +0.20%  " },
    @{ Row = 28; D = $null; E = "  -1.54%  " },
    @{ Row = 29; D = $null; E = "  +0.21%  " },
    @{ Row = 30; D = "'7.95"; E = "  -3.89%  " },
    @{ Row = 31; D = "'2.10"; E = "  -1.92%  " },
    @{ Row = 32; D = "'0.999"; E = "  -0.02%  " },
    @{ Row = 33; D = $null; E = "  -6.20%  " },
    @{ Row = 34; D = "'23.15"; E = "  -2.77%  " },
    @{ Row = 35; D = "'5.31"; E = "  -0.72%  " },
    @{ Row = 36; D = "'1.61"; E = "  +2.08%  " },
    @{ Row = 37; D = "'31.67"; E = "  -0.04%  " },
    @{ Row = 38; D = "'6.95"; E = "  -2.24%  " },
    @{ Row = 39; D = "'170.32"; E = "  -0.52%  " },
    @{ Row = 40; D = "3.471.54"; E = "  -1.54%  " },
    @{ Row = 41; D = "'0.0770"; E = "  +0.03%  " },
    @{ Row = 42; D = "'0.784"; E = "  -3.09%  " },
    @{ Row = 43; D = "'42.51"; E = "  +0.29%  " },
    @{ Row = 44; D = $null; E = "  -1.58%  " },
    @{ Row = 45; D = $null; E = "  -3.28%  " },
    @{ Row = 46; D = "'1.17"; E = "  -3.10%  " },
    @{ Row = 47; D = "2.562.67"; E = "  -2.54%  " },
    @{ Row = 48; D = "'6.91"; E = "  +1.91%  " },
    @{ Row = 49; D = $null; E = "  -1.22%  " },
    @{ Row = 50; D = "'22.56"; E = "  -3.81%  " },
    @{ Row = 51; D = "'0.998"; E = "  -0.32%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $ws.Cells.Item($r, 4).Value = $u.D
        $ws.Cells.Item($r, 4).Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
